$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 1 - "Closed Form for Fibonacci Recurrence" title:
#    merge the two identically-formatted runs ("Closed Form " and
#    "for Fibonacci Recurrence") into a single run.
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Closed Form for Fibonacci Recurrence") {
            # Re-typing the same text collapses the two runs (which already
            # share identical rPr) into a single run.
            $tr.Text = "X"
            $tr.Text = "Closed Form for Fibonacci Recurrence"
        }
    }
}

# ------------------------------------------------------------------
# 2) Recolor the "a"/"b" variable call-outs from orange (FF6600) to
#    red (FF0000) on the two "Closed Form for [xn]B(x)" slides.
# ------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        $full = $shp.TextFrame.TextRange.Text
        if ($full -eq "need to solve for a and b") {
            $tr = $shp.TextFrame.TextRange
            $tr.Characters(19, 1).Font.Color.RGB = 255   # "a"
            $tr.Characters(25, 1).Font.Color.RGB = 255   # "b"
        }
        elseif ($full -eq "Solve for a and b ") {
            $tr = $shp.TextFrame.TextRange
            $tr.Characters(11, 1).Font.Color.RGB = 255   # "a"
            $tr.Characters(17, 1).Font.Color.RGB = 255   # "b"
        }
    }
}

# ------------------------------------------------------------------
# 3) Slide 13 animation timing: the first effect in the main sequence
#    (the "Object 5" equation) switches from "After Previous" to
#    "With Previous".
# ------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$seq = $slide13.TimeLine.MainSequence
for ($i = 1; $i -le $seq.Count; $i++) {
    $eff = $seq.Item($i)
    if ($eff.Timing.TriggerType -eq 3) {
        $eff.Timing.TriggerType = 2
    }
}
